$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: mark "Done", drop bold (keep border) ---
$ws.Range("B6").Value = "Done"
$ws.Range("A6:G6").Font.Bold = $false

# --- Row 7: mark "Done", drop bold (keep border) ---
$ws.Range("B7").Value = "Done"
$ws.Range("A7:G7").Font.Bold = $false

# --- Row 10: new values + drop bold (keep border) ---
$ws.Range("B10").Value = "Canceled"
$ws.Range("F10").Value = "I drop venv support"
$ws.Range("A10:G10").Font.Bold = $false

# --- Restore the selection shown when the file was last saved ---
$ws.Range("E17").Select()
